$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52
$prev = $row - 1

# The new task row repeats the layout/formatting used by the rows above it
# (column A carries the task-name style, column E carries the date style).
$ws.Range("A$prev").Copy() | Out-Null
$ws.Range("A$row").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("E$prev").Copy() | Out-Null
$ws.Range("E$row").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = "Horas entre sesiones"
$ws.Cells.Item($row, 2).Value = 1
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = "Clara"
$ws.Cells.Item($row, 5).Value = (Get-Date -Year 2016 -Month 11 -Day 28).Date

$ws.Range("E$row").Select()
